$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.850.96"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "2.538.78"
$ws.Range("E3").Value = "  -4.60%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.48"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.06"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.50"
$ws.Range("E9").Value = "  -6.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("E10").Value = "  -3.78%  "
$ws.Range("E11").Value = "  -3.81%  "
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "2.985.55"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("D14").Value = "56.869.47"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.02"
$ws.Range("E15").Value = "  -4.90%  "
$ws.Range("E16").Value = "  -3.26%  "
$ws.Range("D17").Value = "2.561.02"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "332.11"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.13"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.25"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.165"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.400"
$ws.Range("E26").Value = "  -4.46%  "
$ws.Range("D27").Value = "2.653.11"
$ws.Range("E27").Value = "  -4.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.87"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "0.0₃0747"
$ws.Range("E29").Value = "  -6.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("E31").Value = "  -6.62%  "
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.45"
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.94"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.838"
$ws.Range("E37").Value = "  -6.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.52"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  -6.09%  "
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0951"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.61"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.577"
$ws.Range("E45").Value = "  -6.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "260.30"
$ws.Range("E46").Value = "  -5.51%  "
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.53"
$ws.Range("E48").Value = "  -6.86%  "
$ws.Range("D49").Value = "1.965.08"
$ws.Range("E49").Value = "  -4.07%  "
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.52"
$ws.Range("E51").Value = "  -3.96%  "
